# Fruta / hortaliza, semanal
# Weekly refresh of daily price data for "Bruselas (repollito)" rows 2-15:
# re-shuffles the Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) values across the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44453
$ws.Range("J2").Value = 280
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 21286
$ws.Range("P2").Value = 1419

$ws.Range("D3").Value = 44398
$ws.Range("J3").Value = 130
$ws.Range("K3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 20000
$ws.Range("P3").Value = 1333

$ws.Range("D4").Value = 44446
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 22000
$ws.Range("L4").Value = 24000
$ws.Range("M4").Value = 22667
$ws.Range("P4").Value = 1511

$ws.Range("D5").Value = 44391
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 20000
$ws.Range("P5").Value = 1333

$ws.Range("D6").Value = 44449
$ws.Range("J6").Value = 220
$ws.Range("K6").Value = 22000
$ws.Range("L6").Value = 24000
$ws.Range("M6").Value = 23091
$ws.Range("P6").Value = 1539

$ws.Range("D7").Value = 44392
$ws.Range("J7").Value = 220
$ws.Range("K7").Value = 23000
$ws.Range("L7").Value = 23000
$ws.Range("M7").Value = 23000
$ws.Range("P7").Value = 1533

$ws.Range("D8").Value = 44400
$ws.Range("J8").Value = 130
$ws.Range("K8").Value = 24000
$ws.Range("L8").Value = 24000
$ws.Range("M8").Value = 24000
$ws.Range("P8").Value = 1600

$ws.Range("D9").Value = 44476
$ws.Range("J9").Value = 220
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 22000
$ws.Range("M9").Value = 20909
$ws.Range("P9").Value = 1394

$ws.Range("D10").Value = 44399
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 22000
$ws.Range("L10").Value = 22000
$ws.Range("M10").Value = 22000
$ws.Range("P10").Value = 1467

$ws.Range("D11").Value = 44483
$ws.Range("J11").Value = 220
$ws.Range("K11").Value = 18000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 18909
$ws.Range("P11").Value = 1261

$ws.Range("D12").Value = 44365
$ws.Range("J12").Value = 580
$ws.Range("K12").Value = 20000
$ws.Range("L12").Value = 22000
$ws.Range("M12").Value = 21103
$ws.Range("P12").Value = 1407

$ws.Range("D13").Value = 44406
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 22000
$ws.Range("M13").Value = 20850
$ws.Range("P13").Value = 1390

$ws.Range("D14").Value = 44396
$ws.Range("J14").Value = 130
$ws.Range("K14").Value = 22000
$ws.Range("L14").Value = 22000
$ws.Range("M14").Value = 22000
$ws.Range("P14").Value = 1467

$ws.Range("D15").Value = 44435
$ws.Range("J15").Value = 140
$ws.Range("K15").Value = 21000
$ws.Range("L15").Value = 23000
$ws.Range("M15").Value = 21714
$ws.Range("P15").Value = 1448
